$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows whose Target cluster is "MuSCs" (rows 12, 9, 6, 3),
# deleting bottom-up so the remaining row numbers stay stable.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(3).Delete()

# Refresh the numeric TPM-derived columns (G:T) for the remaining 8 data rows
# with the newly computed values.
$ws.Range("G2").Value = 3.015833333333333
$ws.Range("H2").Value = 9.047499999999999
$ws.Range("I2").Value = 0.05376901095572644
$ws.Range("J2").Value = 0.05376901095572643
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1797713333333334
$ws.Range("N2").Value = 0.5393140000000001
$ws.Range("O2").Value = 0.188800001120238
$ws.Range("P2").Value = 0.188800001120238
$ws.Range("Q2").Value = 0.5421603794444445
$ws.Range("R2").Value = 4.879443415
$ws.Range("S2").Value = 0.01015158932867524
$ws.Range("T2").Value = 0.01015158932867524
$ws.Range("G3").Value = 3.015833333333333
$ws.Range("H3").Value = 9.047499999999999
$ws.Range("I3").Value = 0.05376901095572644
$ws.Range("J3").Value = 0.05376901095572643
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7724073333333333
$ws.Range("N3").Value = 2.317222
$ws.Range("O3").Value = 0.811199998879762
$ws.Range("P3").Value = 0.811199998879762
$ws.Range("Q3").Value = 2.329451782777777
$ws.Range("R3").Value = 20.965066045
$ws.Range("S3").Value = 0.0436174216270512
$ws.Range("T3").Value = 0.04361742162705119
$ws.Range("G4").Value = 2.006702333333334
$ws.Range("H4").Value = 6.020107
$ws.Range("I4").Value = 0.03577730856453667
$ws.Range("J4").Value = 0.03577730856453666
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1797713333333334
$ws.Range("N4").Value = 0.5393140000000001
$ws.Range("O4").Value = 0.188800001120238
$ws.Range("P4").Value = 0.188800001120238
$ws.Range("Q4").Value = 0.3607475540664445
$ws.Range("R4").Value = 3.246727986598001
$ws.Range("S4").Value = 0.006754755897063622
$ws.Range("T4").Value = 0.006754755897063621
$ws.Range("G5").Value = 2.006702333333334
$ws.Range("H5").Value = 6.020107
$ws.Range("I5").Value = 0.03577730856453667
$ws.Range("J5").Value = 0.03577730856453666
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7724073333333333
$ws.Range("N5").Value = 2.317222
$ws.Range("O5").Value = 0.811199998879762
$ws.Range("P5").Value = 0.811199998879762
$ws.Range("Q5").Value = 1.549991598083778
$ws.Range("R5").Value = 13.949924382754
$ws.Range("S5").Value = 0.02902255266747305
$ws.Range("T5").Value = 0.02902255266747304
$ws.Range("G6").Value = 1.732509666666666
$ws.Range("H6").Value = 5.197528999999999
$ws.Range("I6").Value = 0.03088875310789786
$ws.Range("J6").Value = 0.03088875310789786
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1797713333333334
$ws.Range("N6").Value = 0.5393140000000001
$ws.Range("O6").Value = 0.188800001120238
$ws.Range("P6").Value = 0.188800001120238
$ws.Range("Q6").Value = 0.3114555727895555
$ws.Range("R6").Value = 2.803100155106
$ws.Range("S6").Value = 0.005831796621373868
$ws.Range("T6").Value = 0.005831796621373868
$ws.Range("G7").Value = 1.732509666666666
$ws.Range("H7").Value = 5.197528999999999
$ws.Range("I7").Value = 0.03088875310789786
$ws.Range("J7").Value = 0.03088875310789786
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7724073333333333
$ws.Range("N7").Value = 2.317222
$ws.Range("O7").Value = 0.811199998879762
$ws.Range("P7").Value = 0.811199998879762
$ws.Range("Q7").Value = 1.338203171604222
$ws.Range("R7").Value = 12.043828544438
$ws.Range("S7").Value = 0.02505695648652399
$ws.Range("T7").Value = 0.02505695648652399
$ws.Range("G8").Value = 49.33364366666667
$ws.Range("H8").Value = 148.000931
$ws.Range("I8").Value = 0.879564927371839
$ws.Range("J8").Value = 0.879564927371839
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1797713333333334
$ws.Range("N8").Value = 0.5393140000000001
$ws.Range("O8").Value = 0.188800001120238
$ws.Range("P8").Value = 0.188800001120238
$ws.Range("Q8").Value = 8.868774900148225
$ws.Range("R8").Value = 79.81897410133402
$ws.Range("S8").Value = 0.1660618592731252
$ws.Range("T8").Value = 0.1660618592731252
$ws.Range("G9").Value = 49.33364366666667
$ws.Range("H9").Value = 148.000931
$ws.Range("I9").Value = 0.879564927371839
$ws.Range("J9").Value = 0.879564927371839
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7724073333333333
$ws.Range("N9").Value = 2.317222
$ws.Range("O9").Value = 0.811199998879762
$ws.Range("P9").Value = 0.811199998879762
$ws.Range("Q9").Value = 38.10566814818689
$ws.Range("R9").Value = 342.9510133336821
$ws.Range("S9").Value = 0.7135030680987138
$ws.Range("T9").Value = 0.7135030680987138
